# Append new order-line rows (32-34) to the bottom of the sheet, mirroring
# the context-data rows the updated exporter now emits. Every column in this
# sheet is stored as text (SKU codes, quantities, and prices are all plain
# strings, not numbers), so each value is written with a leading apostrophe
# to force Excel to keep it as literal text instead of auto-converting
# numeric-looking strings like "1" or "0.00" into numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 32: Palmer / Butter - Salted
$ws.Range("A32").Value = "'Palmer"
$ws.Range("B32").Value = "'Butter - Salted"
$ws.Range("C32").Value = "'1"
$ws.Range("D32").Value = "'0.00"
$ws.Range("E32").Value = "'0.00"

# Row 33: PERF / Vegan Egg
$ws.Range("A33").Value = "'PERF"
$ws.Range("B33").Value = "'Vegan Egg"
$ws.Range("C33").Value = "'1"
$ws.Range("D33").Value = "'99.59"
$ws.Range("E33").Value = "'99.59"

# Row 34: (blank SKU) / Flour - Millers Choice
$ws.Range("A34").Value = "'"
$ws.Range("B34").Value = "'Flour - Millers Choice"
$ws.Range("C34").Value = "'1"
$ws.Range("D34").Value = "'0.00"
$ws.Range("E34").Value = "'0.00"
